$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.375839
$ws.Range("H2").Value = 16.127517
$ws.Range("I2").Value = 0.2354568587499626
$ws.Range("J2").Value = 0.2354568587499626
$ws.Range("O2").Value = 0.06522509891308133
$ws.Range("P2").Value = 0.06522509891308133
$ws.Range("Q2").Value = 1.101054256731334
$ws.Range("R2").Value = 9.909488310582001
$ws.Range("S2").Value = 0.01535769690172973
$ws.Range("T2").Value = 0.01535769690172973

# Row 3
$ws.Range("G3").Value = 5.375839
$ws.Range("H3").Value = 16.127517
$ws.Range("I3").Value = 0.2354568587499626
$ws.Range("J3").Value = 0.2354568587499626
$ws.Range("M3").Value = 0.6481333333333333
$ws.Range("N3").Value = 1.9444
$ws.Range("O3").Value = 0.2064033004146749
$ws.Range("P3").Value = 0.2064033004146749
$ws.Range("Q3").Value = 3.484260450533333
$ws.Range("R3").Value = 31.3583440548
$ws.Range("S3").Value = 0.0485990727512642
$ws.Range("T3").Value = 0.0485990727512642

# Row 4
$ws.Range("G4").Value = 5.375839
$ws.Range("H4").Value = 16.127517
$ws.Range("I4").Value = 0.2354568587499626
$ws.Range("J4").Value = 0.2354568587499626
$ws.Range("M4").Value = 1.888205
$ws.Range("N4").Value = 5.664615
$ws.Range("O4").Value = 0.6013141491351952
$ws.Range("P4").Value = 0.6013141491351952
$ws.Range("Q4").Value = 10.150686078995
$ws.Range("R4").Value = 91.356174710955
$ws.Range("S4").Value = 0.1415835406772796
$ws.Range("T4").Value = 0.1415835406772796

# Row 5
$ws.Range("G5").Value = 5.375839
$ws.Range("H5").Value = 16.127517
$ws.Range("I5").Value = 0.2354568587499626
$ws.Range("J5").Value = 0.2354568587499626
$ws.Range("M5").Value = 0.398977
$ws.Range("N5").Value = 1.196931
$ws.Range("O5").Value = 0.1270574515370486
$ws.Range("P5").Value = 0.1270574515370486
$ws.Range("Q5").Value = 2.144836116703
$ws.Range("R5").Value = 19.303525050327
$ws.Range("S5").Value = 0.02991654841968906
$ws.Range("T5").Value = 0.02991654841968906

# Row 6
$ws.Range("I6").Value = 0.007131134316291014
$ws.Range("J6").Value = 0.007131134316291014
$ws.Range("O6").Value = 0.06522509891308133
$ws.Range("P6").Value = 0.06522509891308133
$ws.Range("S6").Value = 0.00046512894114255
$ws.Range("T6").Value = 0.00046512894114255

# Row 7
$ws.Range("I7").Value = 0.007131134316291014
$ws.Range("J7").Value = 0.007131134316291014
$ws.Range("M7").Value = 0.6481333333333333
$ws.Range("N7").Value = 1.9444
$ws.Range("O7").Value = 0.2064033004146749
$ws.Range("P7").Value = 0.2064033004146749
$ws.Range("Q7").Value = 0.1055256126222222
$ws.Range("R7").Value = 0.9497305135999999
$ws.Range("S7").Value = 0.001471889658582812
$ws.Range("T7").Value = 0.001471889658582811

# Row 8
$ws.Range("I8").Value = 0.007131134316291014
$ws.Range("J8").Value = 0.007131134316291014
$ws.Range("M8").Value = 1.888205
$ws.Range("N8").Value = 5.664615
$ws.Range("O8").Value = 0.6013141491351952
$ws.Range("P8").Value = 0.6013141491351952
$ws.Range("Q8").Value = 0.3074274676733333
$ws.Range("R8").Value = 2.76684720906
$ws.Range("S8").Value = 0.004288051963769323
$ws.Range("T8").Value = 0.004288051963769323

# Row 9
$ws.Range("I9").Value = 0.007131134316291014
$ws.Range("J9").Value = 0.007131134316291014
$ws.Range("M9").Value = 0.398977
$ws.Range("N9").Value = 1.196931
$ws.Range("O9").Value = 0.1270574515370486
$ws.Range("P9").Value = 0.1270574515370486
$ws.Range("Q9").Value = 0.06495930726266666
$ws.Range("R9").Value = 0.584633765364
$ws.Range("S9").Value = 0.0009060637527963294
$ws.Range("T9").Value = 0.0009060637527963294

# Row 10
$ws.Range("G10").Value = 9.994147
$ws.Range("H10").Value = 29.982441
$ws.Range("I10").Value = 0.4377345486919088
$ws.Range("J10").Value = 0.4377345486919088
$ws.Range("O10").Value = 0.06522509891308133
$ws.Range("P10").Value = 0.06522509891308133
$ws.Range("Q10").Value = 2.046954549187333
$ws.Range("R10").Value = 18.422590942686
$ws.Range("S10").Value = 0.02855127923610276
$ws.Range("T10").Value = 0.02855127923610276

# Row 11
$ws.Range("G11").Value = 9.994147
$ws.Range("H11").Value = 29.982441
$ws.Range("I11").Value = 0.4377345486919088
$ws.Range("J11").Value = 0.4377345486919088
$ws.Range("M11").Value = 0.6481333333333333
$ws.Range("N11").Value = 1.9444
$ws.Range("O11").Value = 0.2064033004146749
$ws.Range("P11").Value = 0.2064033004146749
$ws.Range("Q11").Value = 6.477539808933333
$ws.Range("R11").Value = 58.29785828039999
$ws.Range("S11").Value = 0.09034985555553819
$ws.Range("T11").Value = 0.09034985555553818

# Row 12
$ws.Range("G12").Value = 9.994147
$ws.Range("H12").Value = 29.982441
$ws.Range("I12").Value = 0.4377345486919088
$ws.Range("J12").Value = 0.4377345486919088
$ws.Range("M12").Value = 1.888205
$ws.Range("N12").Value = 5.664615
$ws.Range("O12").Value = 0.6013141491351952
$ws.Range("P12").Value = 0.6013141491351952
$ws.Range("Q12").Value = 18.870998336135
$ws.Range("R12").Value = 169.838985025215
$ws.Range("S12").Value = 0.2632159776937538
$ws.Range("T12").Value = 0.2632159776937538

# Row 13
$ws.Range("G13").Value = 9.994147
$ws.Range("H13").Value = 29.982441
$ws.Range("I13").Value = 0.4377345486919088
$ws.Range("J13").Value = 0.4377345486919088
$ws.Range("M13").Value = 0.398977
$ws.Range("N13").Value = 1.196931
$ws.Range("O13").Value = 0.1270574515370486
$ws.Range("P13").Value = 0.1270574515370486
$ws.Range("Q13").Value = 3.987434787619
$ws.Range("R13").Value = 35.886913088571
$ws.Range("S13").Value = 0.05561743620651402
$ws.Range("T13").Value = 0.05561743620651402

# Row 14
$ws.Range("G14").Value = 0.7761303333333333
$ws.Range("H14").Value = 2.328391
$ws.Range("I14").Value = 0.03399380269149206
$ws.Range("J14").Value = 0.03399380269149207
$ws.Range("O14").Value = 0.06522509891308133
$ws.Range("P14").Value = 0.06522509891308133
$ws.Range("Q14").Value = 0.1589633929317778
$ws.Range("R14").Value = 1.430670536386
$ws.Range("S14").Value = 0.00221724914298434
$ws.Range("T14").Value = 0.002217249142984341

# Row 15
$ws.Range("G15").Value = 0.7761303333333333
$ws.Range("H15").Value = 2.328391
$ws.Range("I15").Value = 0.03399380269149206
$ws.Range("J15").Value = 0.03399380269149207
$ws.Range("M15").Value = 0.6481333333333333
$ws.Range("N15").Value = 1.9444
$ws.Range("O15").Value = 0.2064033004146749
$ws.Range("P15").Value = 0.2064033004146749
$ws.Range("Q15").Value = 0.5030359400444444
$ws.Range("R15").Value = 4.5273234604
$ws.Range("S15").Value = 0.007016433069169221
$ws.Range("T15").Value = 0.007016433069169221

# Row 16
$ws.Range("G16").Value = 0.7761303333333333
$ws.Range("H16").Value = 2.328391
$ws.Range("I16").Value = 0.03399380269149206
$ws.Range("J16").Value = 0.03399380269149207
$ws.Range("M16").Value = 1.888205
$ws.Range("N16").Value = 5.664615
$ws.Range("O16").Value = 0.6013141491351952
$ws.Range("P16").Value = 0.6013141491351952
$ws.Range("Q16").Value = 1.465493176051667
$ws.Range("R16").Value = 13.189438584465
$ws.Range("S16").Value = 0.02044095454130426
$ws.Range("T16").Value = 0.02044095454130426

# Row 17
$ws.Range("G17").Value = 0.7761303333333333
$ws.Range("H17").Value = 2.328391
$ws.Range("I17").Value = 0.03399380269149206
$ws.Range("J17").Value = 0.03399380269149207
$ws.Range("M17").Value = 0.398977
$ws.Range("N17").Value = 1.196931
$ws.Range("O17").Value = 0.1270574515370486
$ws.Range("P17").Value = 0.1270574515370486
$ws.Range("Q17").Value = 0.3096581520023333
$ws.Range("R17").Value = 2.786923368021
$ws.Range("S17").Value = 0.004319165938034244
$ws.Range("T17").Value = 0.004319165938034245

# Row 18
$ws.Range("G18").Value = 6.522593333333333
$ws.Range("H18").Value = 19.56778
$ws.Range("I18").Value = 0.2856836555503455
$ws.Range("J18").Value = 0.2856836555503455
$ws.Range("O18").Value = 0.06522509891308133
$ws.Range("P18").Value = 0.06522509891308133
$ws.Range("Q18").Value = 1.335927127764444
$ws.Range("R18").Value = 12.02334414988
$ws.Range("S18").Value = 0.01863374469112194
$ws.Range("T18").Value = 0.01863374469112194

# Row 19
$ws.Range("G19").Value = 6.522593333333333
$ws.Range("H19").Value = 19.56778
$ws.Range("I19").Value = 0.2856836555503455
$ws.Range("J19").Value = 0.2856836555503455
$ws.Range("M19").Value = 0.6481333333333333
$ws.Range("N19").Value = 1.9444
$ws.Range("O19").Value = 0.2064033004146749
$ws.Range("P19").Value = 0.2064033004146749
$ws.Range("Q19").Value = 4.227510159111111
$ws.Range("R19").Value = 38.047591432
$ws.Range("S19").Value = 0.05896604938012048
$ws.Range("T19").Value = 0.05896604938012047

# Row 20
$ws.Range("G20").Value = 6.522593333333333
$ws.Range("H20").Value = 19.56778
$ws.Range("I20").Value = 0.2856836555503455
$ws.Range("J20").Value = 0.2856836555503455
$ws.Range("M20").Value = 1.888205
$ws.Range("N20").Value = 5.664615
$ws.Range("O20").Value = 0.6013141491351952
$ws.Range("P20").Value = 0.6013141491351952
$ws.Range("Q20").Value = 12.31599334496667
$ws.Range("R20").Value = 110.8439401047
$ws.Range("S20").Value = 0.1717856242590882
$ws.Range("T20").Value = 0.1717856242590882

# Row 21
$ws.Range("G21").Value = 6.522593333333333
$ws.Range("H21").Value = 19.56778
$ws.Range("I21").Value = 0.2856836555503455
$ws.Range("J21").Value = 0.2856836555503455
$ws.Range("M21").Value = 0.398977
$ws.Range("N21").Value = 1.196931
$ws.Range("O21").Value = 0.1270574515370486
$ws.Range("P21").Value = 0.1270574515370486
$ws.Range("Q21").Value = 2.602364720353333
$ws.Range("R21").Value = 23.42128248318
$ws.Range("S21").Value = 0.0362982372200149
$ws.Range("T21").Value = 0.0362982372200149
